# Auto-save via app Streamlit
# A new reservation ("Annick" / "Autre") was inserted as row 37 of the
# reservations sheet, pushing the existing rows 37-51 (and the TOTAL row)
# down by one (new rows 38-52). This recreates that edit using Excel's
# native row-insert behaviour so every following row (and the trailing
# TOTAL row) shifts down automatically, then fills in the new row's data.
# Along the way, the phone number for "Gerlinde Weiss" (now row 49) is
# normalized from the text "4369912047111.0" to the clean numeric value
# 4369912047111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 37; everything below (including the
# TOTAL row) shifts down by one automatically.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new reservation's data.
$ws.Range("A37").Value = "Annick"
$ws.Range("B37").Value = "Autre"
$ws.Range("D37").Value = 45880
$ws.Range("E37").Value = 45882
$ws.Range("F37").Value = 2
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2025
$ws.Range("L37").Value = 8

# Normalize the Gerlinde Weiss phone number (now on row 49 after the
# shift) to a clean numeric value instead of the text "4369912047111.0".
$ws.Range("C49").Value = 4369912047111
